# Fruta / hortaliza, semanal
#
# Inserts a new weekly price-report row for "Agrícola del Norte S.A. de
# Arica - Plátano" just above the existing row 231, pushing the previously
# recorded rows 231-261 down to 232-262 (Excel's native "insert row"
# behaviour - every subsequent row keeps its own data, it just moves down
# one position). The sheet's used range grows from A1:T261 to A1:T262.
#
# The brand-new row reuses the static/categorical fields that the old
# row 231 already had (market, region, product hierarchy, variety,
# quality, volume, unit, origin, kg/unit) and only carries fresh
# date/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 231..261 down to 232..262, leaving a blank row 231 behind.
$ws.Rows(231).Insert()

# Populate the newly inserted row 231 with the latest weekly report.
$ws.Range("A231").Value = 1
$ws.Range("B231").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C231").Value = "Arica y Parinacota"
$ws.Range("D231").Value = 44748
$ws.Range("E231").Value = 15
$ws.Range("F231").Value = "Fruta"
$ws.Range("G231").Value = 100108
$ws.Range("H231").Value = "Tropicales y subtropicales"
$ws.Range("I231").Value = 100108006
$ws.Range("J231").Value = "Plátano"
$ws.Range("K231").Value = "Sin especificar"
$ws.Range("L231").Value = "Pintón"
$ws.Range("M231").Value = 120
$ws.Range("N231").Value = 25000
$ws.Range("O231").Value = 26000
$ws.Range("P231").Value = 25500
$ws.Range("Q231").Value = "$/caja 20 kilos"
$ws.Range("R231").Value = "Ecuador"
$ws.Range("S231").Value = 1275
$ws.Range("T231").Value = 20
